$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.447.83'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.77%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.218.94'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.18%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.11'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.06%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '182.75'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.89%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.601'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.24%  '

$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.129'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.79%  '

$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.59'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.45%  '

$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.412'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.06%  '

$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.780.86'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.13%  '

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.137'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.15%  '

$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.56'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.49%  '

$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '67.500.47'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.69%  '

$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000168'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.57%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.221.46'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.17%  '

$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.73'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.26%  '

$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.38'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.67%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '391.83'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.49%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.52'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.29%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.11%  '

$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.63'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.09%  '

$ws.Range('B24').Value = 'Polygon'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.511'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.74%  '

$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000117'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.42%  '

$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.186'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.70%  '

$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.50'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.62%  '

$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.15%  '

$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.95'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.26%  '

$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.55'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.74%  '

$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.54'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.73%  '

$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.94'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -5.43%  '

$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.05%  '

$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.24'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.61%  '

$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '161.38'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.69%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.47'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.88%  '

$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.88'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.25%  '

$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '26.18'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.88%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.801'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.05%  '

$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.53'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.77%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.43'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.39%  '

$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.45'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.88%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.70'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.63%  '

$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0679'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.39%  '

$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.596.79'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.82%  '

$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '24.68'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.33%  '

$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '333.84'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.87%  '

$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0276'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.31%  '

$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.25'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.66%  '

$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.101'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.25%  '

$ws.Range('B51').Value = 'Arweave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '30.75'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.00%  '
